# daily auto push: 2026-01-16 09:37 UTC
# Insert a new daily-stats row just before the "2026/12/29" block (current
# row 650), pushing that row and everything after it down by one. The new
# row carries the 2026/01/16 (Fri) reading that was missing from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 650 (and everything below it) down by one row.
$ws.Rows.Item(650).Insert()

# Column A/B hold date-like / kanji text (t="inlineStr" in the source file).
# Excel's normal Value-assignment auto-converts a "yyyy/mm/dd"-shaped string
# into a real date serial, so force Text entry via NumberFormat, then strip
# the format back off (ClearFormats) so the cell is left with no explicit
# style -- matching every other data row on the sheet.
$ws.Range("A650").NumberFormat = "@"
$ws.Range("A650").Value = "2026/01/16"
$ws.Range("A650").ClearFormats()

$ws.Range("B650").NumberFormat = "@"
$ws.Range("B650").Value = "金"
$ws.Range("B650").ClearFormats()

$ws.Range("C650").Value = 17
$ws.Range("D650").Value = 201
